$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E3").Value = 16.697
$ws.Range("B4").Value = 7.723999999999999
$ws.Range("B6").Value = 6.778
$ws.Range("B7").Value = 5.491999999999999
$ws.Range("C7").Value = -13.218
$ws.Range("B8").Value = 6.238999999999999
$ws.Range("C11").Value = -12.767
$ws.Range("C12").Value = -11.929
$ws.Range("D12").Value = -7.373
$ws.Range("E12").Value = 17.233
$ws.Range("D13").Value = -7.63
$ws.Range("E13").Value = 16.682
$ws.Range("D14").Value = -7.781999999999999
$ws.Range("C15").Value = -13.102
$ws.Range("B16").Value = 5.817
$ws.Range("D16").Value = -8.171999999999999
$ws.Range("D19").Value = -7.944000000000001
$ws.Range("B20").Value = 9.244
$ws.Range("C20").Value = -12.165
$ws.Range("D20").Value = -7.962000000000001
$ws.Range("B21").Value = 9.205
$ws.Range("C21").Value = -12.12
$ws.Range("C22").Value = -12.912
$ws.Range("D22").Value = -7.375
$ws.Range("E22").Value = 16.896
$ws.Range("C23").Value = -12.731
$ws.Range("E25").Value = 17.179
$ws.Range("B28").Value = 5.881
$ws.Range("B29").Value = 5.406000000000001
$ws.Range("C29").Value = -11.641
$ws.Range("E29").Value = 17.286
$ws.Range("B30").Value = 5.300000000000001
$ws.Range("B32").Value = 5.918999999999999
$ws.Range("C34").Value = -12.974
$ws.Range("E34").Value = 16.823
$ws.Range("D36").Value = -7.855
$ws.Range("B40").Value = 9.203999999999999
$ws.Range("C42").Value = -12.092
$ws.Range("C43").Value = -13.704
$ws.Range("D43").Value = -7.893999999999998
$ws.Range("E43").Value = 16.597
$ws.Range("C44").Value = -13.339
$ws.Range("C45").Value = -13.262
$ws.Range("B46").Value = 5.545
$ws.Range("C46").Value = -14.006
$ws.Range("D46").Value = -8.562000000000001
$ws.Range("E48").Value = 17.022
$ws.Range("C50").Value = -13.003
$ws.Range("D50").Value = -8.518000000000001
$ws.Range("B51").Value = 5.447000000000001
$ws.Range("C51").Value = -12.454
$ws.Range("B52").Value = 5.82
$ws.Range("B57").Value = 6.358
$ws.Range("C57").Value = -14.252
$ws.Range("B59").Value = 6.098999999999999
$ws.Range("E60").Value = 16.403
$ws.Range("B62").Value = 6.358
$ws.Range("C65").Value = -12.45
$ws.Range("B66").Value = 4.998
$ws.Range("C66").Value = -10.868
$ws.Range("C67").Value = -11.926
$ws.Range("E68").Value = 17.337
$ws.Range("E70").Value = 17.46
$ws.Range("E71").Value = 17.117
$ws.Range("B73").Value = 7.521000000000001
$ws.Range("E73").Value = 16.693
$ws.Range("B74").Value = 9.132999999999999
$ws.Range("D76").Value = -7.383000000000001
$ws.Range("B77").Value = 6.641
$ws.Range("E78").Value = 16.931
$ws.Range("C79").Value = -12.863
$ws.Range("C84").Value = -13.585
$ws.Range("C87").Value = -13.479
$ws.Range("E87").Value = 16.17
$ws.Range("B92").Value = 5.959000000000001
$ws.Range("C92").Value = -12.628
$ws.Range("E92").Value = 17.271
$ws.Range("D95").Value = -7.536
$ws.Range("C97").Value = -12.291
$ws.Range("D97").Value = -7.795999999999999
$ws.Range("D99").Value = -7.520999999999999
$ws.Range("B100").Value = 6.651999999999999
$ws.Range("E101").Value = 16.737
